$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E4"   = 16.37599999999999
    "E7"   = 15.65390000000001
    "E16"  = 16.0584
    "E28"  = 16.51029999999999
    "E29"  = 17.04380000000002
    "E32"  = 16.94979999999999
    "E40"  = 17.06270000000001
    "E52"  = 16.99990000000001
    "E57"  = 16.62169999999999
    "E66"  = 17.19970000000002
    "E100" = 16.3964
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
